$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O1").EntireColumn.Insert()
$ws.Range("O1").Value = "literal"
